# Update cryptocurrency price (column D) and volume/change (column E) values
# to reflect the latest scrape, per GitHub Actions automation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "69.292.66"
$d.Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "3.741.99"
$d.Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  -0.10%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "614.42"
$d.Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = "178.46"
$d.Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "3.739.59"
$d.Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -1.85%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "0.166"
$d.Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "6.55"
$d.Style = "Normal"
$ws.Range("E11").Value = "  +3.53%  "

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "0.480"
$d.Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "39.88"
$d.Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "0.0000253"
$d.Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "4.361.55"
$d.Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "3.736.31"
$d.Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "69.366.21"
$d.Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("E18").Value = "  -2.45%  "

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "7.41"
$d.Style = "Normal"
$ws.Range("E19").Value = "  -1.76%  "

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "16.33"
$d.Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "498.64"
$d.Style = "Normal"
$ws.Range("E21").Value = "  -2.62%  "

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "9.18"
$d.Style = "Normal"
$ws.Range("E22").Value = "  -3.22%  "

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "0.720"
$d.Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "2.57"
$d.Style = "Normal"
$ws.Range("E24").Value = "  +3.12%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "85.56"
$d.Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "12.89"
$d.Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "10.92"
$d.Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "0.0000134"
$d.Style = "Normal"
$ws.Range("E28").Value = "  +6.40%  "

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "1.00"
$d.Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "2.49"
$d.Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = "2.91"
$d.Style = "Normal"
$ws.Range("E31").Value = "  +2.81%  "

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "8.02"
$d.Style = "Normal"
$ws.Range("E32").Value = "  +3.03%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "30.31"
$d.Style = "Normal"
$ws.Range("E33").Value = "  -3.17%  "

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "0.113"
$d.Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "0.998"
$d.Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("E36").Value = "  +0.50%  "

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "6.11"
$d.Style = "Normal"
$ws.Range("E37").Value = "  -1.27%  "

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "0.348"
$d.Style = "Normal"
$ws.Range("E38").Value = "  +3.25%  "

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "0.136"
$d.Style = "Normal"
$ws.Range("E39").Value = "  +3.20%  "

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = "3.08"
$d.Style = "Normal"
$ws.Range("E40").Value = "  +13.43%  "

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "447.47"
$d.Style = "Normal"
$ws.Range("E41").Value = "  +6.88%  "

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = "2.06"
$d.Style = "Normal"
$ws.Range("E42").Value = "  -4.79%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "49.73"
$d.Style = "Normal"
$ws.Range("E43").Value = "  -2.73%  "

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "44.51"
$d.Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "8.54"
$d.Style = "Normal"
$ws.Range("E45").Value = "  -2.47%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "2.945.80"
$d.Style = "Normal"
$ws.Range("E46").Value = "  -4.07%  "

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "0.0359"
$d.Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "

$ws.Range("E48").Value = "  -0.01%  "

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "138.41"
$d.Style = "Normal"
$ws.Range("E49").Value = "  +1.99%  "

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "27.03"
$d.Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "

$ws.Range("E51").Value = "  -1.61%  "
